$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("queries")

# New header columns for the "queries" sheet: param.format, param.q
$ws.Cells.Item(1, 4).Value = "param.format"
$ws.Cells.Item(1, 5).Value = "param.q"

# New row 5: odk_values query example (row 4 intentionally left blank)
$ws.Cells.Item(5, 1).Value = "odk_values"
$ws.Cells.Item(5, 2).Value = "odkquery://table_id/elementKey1/elementKey5/?selection=encodeURIComponent('elementKey2=? and elementKey3>5')&selectionArgs=encodeURIComponent(JSON.stringify([data('state')])"
